# Appendix 4 - Data for simulations: add human_population_2024,
# gdp_per_capita_20224 and gdp_2024 columns to the Pop_Data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pop_Data")

# New column headers (row 1)
$ws.Range("V1").Value = "human_population_2024"
$ws.Range("W1").Value = "gdp_per_capita_20224"
$ws.Range("X1").Value = "gdp_2024"

$v_vals = @(46731365.000000015,37298625,13470099.999999998,2525764.9999999995,23351513.000000019,13378732.000000002,29314689.000000004,5226574.0000000019,18437024.999999989,101115264.99999997,6086280.0000000009,29159502.999999996,1043016.0000000001,109887164,1590203.9999999995,3795984.0000000009,1210452.0000000002,126750867.99999991,2435945.0000000009,2705194.0000000005,33733902.000000007,14617783.999999994,2160279.9999999991,58714808.999999978,2209599.0000000023,5561376.9999999981,7193917,21287587.999999996,22752839.999999996,5159892.9999999953,38606160,34999755.999999978,2727257,28079834.000000004,227713026.99999997,14249748.000000006,18585815,8638684.0000000019,17851254.999999996,62134218.000000037,12176470.999999998,48216428.999999978,67014130.000000045,9095785.0000000075,12251377.999999994,50976034.000000015,20608355.999999996,15840801.999999994)

$w_vals = @(4114.1549906144819,2336.1513177117904,1560.9027980683363,8028.8688720748514,1003.2901163852914,258.75274253438215,1815.7867563857214,558.91319159676959,760.99997023889375,638.26884422319938,2419.1673936541579,2817.8811287763206,3675.6189752318733,4235.802789781611,9246.96875410752,689.04304984399994,4605.6931675416017,1031.4967902263556,8760.3819731424155,913.07658530156198,2672.0352246405087,1283.2471999743909,888.34593496511081,2192.9197732123466,1274.6241575173519,735.50259700305287,6576.517576235141,702.24856361256707,1003.0283878665916,1882.7837724833846,3821.002207737793,546.84955428605815,5167.8020435604667,650.09072944156355,2278.3694125050547,911.1484185545238,1755.4353956616928,563.77292691362845,487.1147673341577,7642.7639175561035,398.17879153000007,835.21605164190703,1240.8348600960667,1084.3440676319865,4288.2365639374202,937.62261171112254,1224.54284434702,1898.2569083438891)

$x_vals = @(192260078532.97699,87135231942.587921,21025516780.260296,20279035986.676132,23428342195.542664,3461783596.6325002,53229224053.766197,2921201155.4566956,14030575476.293732,64538723324.872505,14723730124.649429,82168013228.196503,3833729401.0704479,465460355832.3894,14704566700.65679,2615596392.5190268,5574970506.0370684,130743113500.4044,21339808665.56641,2470049300.0982742,90138174408.570923,18758230387.830444,1919075956.3864288,128756865636.48621,2816408263.8261862,4090407226.4130459,47310921592.476776,14949178095.776117,22821744424.586498,9714962808.1506004,147514222592.27847,19139600968.720776,14093924297.914587,18254439767.65802,518814395545.73761,12983635355.000494,32626197508.220024,4870256163.3619328,8695609925.9477177,474877159375.96527,4848412507.8800907,40271135453.652328,83153468623.009674,9862960505.2060165,52536807098.218483,47796282133.754997,25235814873.555973,30069911830.207684)

for ($i = 0; $i -lt 48; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 22).Value = $v_vals[$i]
    $ws.Cells.Item($row, 23).Value = $w_vals[$i]
    $ws.Cells.Item($row, 24).Value = $x_vals[$i]
}

# Column V (human_population_2024) carries Excel's built-in "Comma [0]" cell
# style (thousands separator, no decimals).
$ws.Range("V1:V49").Style = "Comma [0]"

# Column X (gdp_2024) gets the same visual thousands-separator formatting but
# as a direct number format (no named cell style attached). The header cell
# (X1) stays in the sheet's default/general format.
$commaFmt = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"
$ws.Range("X2:X49").NumberFormat = $commaFmt

# Autofit the new columns to mirror the bestFit widths Excel would apply
$ws.Columns.Item(22).AutoFit() | Out-Null
$ws.Columns.Item(23).AutoFit() | Out-Null
$ws.Columns.Item(24).AutoFit() | Out-Null

# Update selection / view to mirror final state in the source workbook
$ws.Range("X1:X1048576").Select()
